# feat: add 2022-Q4 data
#
# - Insert a new "2022-Q4" worksheet right after "总计", holding the latest
#   quarter's fund-holding snapshot (built as a copy of the "2022-Q2" sheet
#   so it inherits the exact same header/column styling).
# - Update the "总计" (overview) sheet: insert a new top data row for
#   2022-Q4 and push the existing 2022-Q2 / 2021-Q2 / 2020-Q4 rows down by
#   one, keeping their original figures intact.
# - The pre-existing "2022-Q2", "2021-Q2", "2020-Q4" sheets keep their data
#   untouched; they simply shift one tab to the right to make room.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Build the new "2022-Q4" sheet from a copy of "2022-Q2" (same layout
#    and cell styles), then rename + relocate it right after "总计".
# ---------------------------------------------------------------------
$sourceSheet = $wb.Worksheets.Item("2022-Q2")
$sourceSheet.Copy([System.Reflection.Missing]::Value, $sourceSheet)
$q4SheetTemp = $wb.Worksheets.Item("2022-Q2 (2)")
$q4SheetTemp.Name = "2022-Q4"

$zongji = $wb.Worksheets.Item("总计")
$q4SheetTemp.Move([System.Reflection.Missing]::Value, $zongji)

# NOTE: re-fetch the sheet handle by name after Move() -- this COM bridge
# rebinds previously-held worksheet references to their old positional
# index rather than following the moved sheet, so reusing $q4SheetTemp
# here would silently edit whatever sheet now sits in its old slot.
$q4Sheet = $wb.Worksheets.Item("2022-Q4")

# Fill in the 2022-Q4 fund-holding figures (A2/A3 already carry the right
# index values + styling from the copied sheet, so only B..H need to change).
# Force text storage for the numeric-looking fund codes / percentages (B, D-G)
# so leading zeros and fixed decimal formatting survive, matching the source data.
$q4Sheet.Range("B2:B3").NumberFormat = "@"
$q4Sheet.Range("D2:G3").NumberFormat = "@"

$q4Sheet.Range("B2").Value = "011205"
$q4Sheet.Range("C2").Value = "兴银中证500指数增强C"
$q4Sheet.Range("D2").Value = "0.66"
$q4Sheet.Range("E2").Value = "84.84"
$q4Sheet.Range("F2").Value = "0.79"
$q4Sheet.Range("G2").Value = "0.0052"
$q4Sheet.Range("H2").Value = 3

$q4Sheet.Range("B3").Value = "010253"
$q4Sheet.Range("C3").Value = "兴银中证500指数增强A"
$q4Sheet.Range("D3").Value = "0.32"
$q4Sheet.Range("E3").Value = "84.84"
$q4Sheet.Range("F3").Value = "0.79"
$q4Sheet.Range("G3").Value = "0.0025"
$q4Sheet.Range("H3").Value = 3

# ---------------------------------------------------------------------
# 2) Update the "总计" overview sheet: insert a fresh row 2 for 2022-Q4,
#    pushing the old rows (2022-Q2, 2021-Q2, 2020-Q4) down to rows 3-5.
# ---------------------------------------------------------------------
$zongji.Rows.Item(2).Insert()

# The insert leaves row 2 without the data rows' normal styling; copy it
# over from row 3 (still carrying the original per-column style) first.
$zongji.Range("A3:D3").Copy()
$zongji.Range("A2:D2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Write every data row's final values explicitly rather than relying on the
# row-shift to carry the right index along (the source data keeps the 0/1/2
# index sequence tied to content, not to physical row, with a fresh "3"
# appended for the newly-created last row).
$zongji.Range("A2").Value = 0
$zongji.Range("B2").Value = "2022-Q4"
$zongji.Range("C2").Value = 2
$zongji.Range("D2").Value = 0.01

$zongji.Range("A3").Value = 1
$zongji.Range("B3").Value = "2022-Q2"
$zongji.Range("C3").Value = 2
$zongji.Range("D3").Value = 0.16

$zongji.Range("A4").Value = 2
$zongji.Range("B4").Value = "2021-Q2"
$zongji.Range("C4").Value = 2
$zongji.Range("D4").Value = 0.06

$zongji.Range("A5").Value = 3
$zongji.Range("B5").Value = "2020-Q4"
$zongji.Range("C5").Value = 5
$zongji.Range("D5").Value = 0.02
